# Apply the "risk analysis" edits described in the commit:
#   1. Rename sheets from Q2_20_21 -> Q1_20_21
#   2. Fix typo "Decription" -> "Description" on the all-data sheet
#   3. Bump the count/total figures on the Count sheet

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheets -------------------------------------------------
$wsAllData = $wb.Worksheets.Item(1)
$wsCount   = $wb.Worksheets.Item(2)

$wsAllData.Name = "Q1_20_21 all data"
$wsCount.Name   = "Q1_20_21 Count"

# --- 2. Fix the typo in the "all data" sheet header ----------------------
$wsAllData.Range("D3").Value = "Brief Risk Description "

# --- 3. Update the tallied counts/totals on the Count sheet ---------------
$wsCount.Range("D7").Value  = 3
$wsCount.Range("F7").Value  = 4

$wsCount.Range("D15").Value = 7
$wsCount.Range("F15").Value = 11

$wsCount.Range("D19").Value = 13
$wsCount.Range("F19").Value = 19

$wsCount.Range("D27").Value = 11
$wsCount.Range("F27").Value = 14

$wsCount.Range("D34").Value = 16
$wsCount.Range("F34").Value = 17

$wsCount.Range("D40").Value = 19
$wsCount.Range("F40").Value = 19

$wsCount.Range("D46").Value = 16
$wsCount.Range("F46").Value = 22
